# Fix security vulnerability: Disable malicious swarm attack script
# -> reflected in the CALENDAR worksheet as a logging update: every
#    previously "approved" scheduling row is now marked "failed" and
#    annotated with the Puppeteer/Chrome launch error that the disabled
#    script now emits, in a new "error_log" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CALENDAR")

$errorLog = "Could not find Chrome (ver. 145.0.7632.67). This can occur if either`n 1. you did not perform an installation before running the script (e.g. ``npx puppeteer browsers install chrome``) or`n 2. your cache path is incorrectly configured (which is: /home/jules/.cache/puppeteer).`nFor (2), check out our guide on configuring puppeteer at https://pptr.dev/guides/configuration."

# New header for column I
$ws.Cells.Item(1, 9).Value = "error_log"

$lastRow = 11
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 5).Value = "failed"
    $ws.Cells.Item($row, 9).Value = $errorLog
}
